# Add a new paragraph "Mailo se despide parte4" after the "Mailo parte3"
# paragraph, mirroring that paragraph's structure (a proofErr-wrapped
# "Mailo" run followed by the rest of the sentence in a second run).

$d = $word.ActiveDocument

# The "Mailo ... parte3" paragraph is currently the last paragraph in the
# document body; insert a brand-new paragraph right after it.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# That insertion created a new (empty) paragraph at the end of the
# document - grab it and fill it in with the exact run/proofErr structure
# from the diff via InsertXML, which replaces the target range's content
# with the supplied WordprocessingML fragment.
$newPara = $d.Paragraphs.Last
$newPara.Range.InsertXML(
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Mailo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> se despide parte4</w:t></w:r>' +
    '</w:p>'
)
